# Adds two new columns, I ("I0") and J ("IF"), to Sheet1, mirroring the
# existing data columns (header in row 1, values in rows 2-63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered/top alignment)
# from the existing H1 header cell so the new headers match the rest of
# the header row exactly (reuses the same cell style).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data ----------------------------------------------------------------
$iValues = @(9,9,9,8,8,9,8,9,8,7,9,12,7,8,7,11,9,9,8,8,8,8,7,8,8,8,9,8,8,9,8,8,8,8,8,7,7,7,7,7,8,7,8,8,8,8,6,7,9,9,8,8,8,7,7,7,8,7,8,7,7,7)
$jValues = @(10,9,9,8,8,9,9,9,8,7,10,12,8,8,7,11,9,9,8,8,8,8,7,8,8,8,9,8,8,9,8,8,8,8,8,7,7,7,8,7,8,7,8,8,8,8,7,7,9,9,8,8,8,7,8,8,8,8,8,7,7,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

Write-Host "Added columns I (I0) and J (IF) for rows 1-63"
